$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 03:52"

# Estados Unidos (row 4) - updated case counts
$ws.Range("B4").Value = 560433
$ws.Range("C4").Value = 133
$ws.Range("E4").Value = 505684
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 22115

# Corea del Sur (row 22) - updated case counts
$ws.Range("B22").Value = 10537
$ws.Range("C22").Value = 25
$ws.Range("D22").Value = 7447
$ws.Range("E22").Value = 2873
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 217

# Australia (row 32) - updated case counts
$ws.Range("B32").Value = 6322
$ws.Range("C32").Value = 9
$ws.Range("E32").Value = 2925

# Mexico's updated totals push it above Filipinas, Arabia Saudita and
# Indonesia in the ranking (sorted descending by total cases), so those
# three rows shift down by one and Mexico takes row 38.
$ws.Range("A38").Value = "Mexico"
$ws.Range("B38").Value = 4661
$ws.Range("C38").Value = 442
$ws.Range("D38").Value = 1772
$ws.Range("E38").Value = 2593
$ws.Range("F38").Value = 89
$ws.Range("G38").Value = 23
$ws.Range("H38").Value = 296

$ws.Range("A39").Value = "Filipinas"
$ws.Range("B39").Value = 4648
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 197
$ws.Range("E39").Value = 4154
$ws.Range("F39").Value = 1
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 297

$ws.Range("A40").Value = "Arabia Saudita"
$ws.Range("B40").Value = 4462
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 761
$ws.Range("E40").Value = 3642
$ws.Range("F40").Value = 67
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 59

$ws.Range("A41").Value = "Indonesia"
$ws.Range("B41").Value = 4241
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 359
$ws.Range("E41").Value = 3509
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 373

# Nueva Zelanda (row 65) - updated case counts
$ws.Range("B65").Value = 1349
$ws.Range("C65").Value = 19
$ws.Range("D65").Value = 546
$ws.Range("E65").Value = 798
$ws.Range("F65").Value = 4
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 5
